# Toggle the "Reviewed Passed" (H) and "Reviewed Bulleted" (I) boolean
# values for the specific rows that were re-reviewed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5, 6, 7, 11, 16, 17, 18, 21, 22, 25)

foreach ($r in $rows) {
    $hCell = $ws.Cells.Item($r, 8)   # Column H
    $iCell = $ws.Cells.Item($r, 9)   # Column I

    $hVal = $hCell.Value2
    $iVal = $iCell.Value2

    $hCell.Value2 = $iVal
    $iCell.Value2 = $hVal
}
